$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024845523301198
$ws.Range("D2").Value = 1.028843176377545
$ws.Range("E2").Value = 1.04939128669869
$ws.Range("F2").Value = 1.0539553926095
$ws.Range("I2").Value = 1.031261377530222
$ws.Range("J2").Value = 1.03001785757573
$ws.Range("K2").Value = 1.031658815516127
$ws.Range("L2").Value = 1.052148475415983
$ws.Range("M2").Value = 1.056699933778067
$ws.Range("N2").Value = 1.013997919318864
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.025935959692008
$ws.Range("D3").Value = 1.029617139707675
$ws.Range("E3").Value = 1.050543389608113
$ws.Range("F3").Value = 1.055150252087689
$ws.Range("I3").Value = 1.031449767668057
$ws.Range("J3").Value = 1.030746646645621
$ws.Range("K3").Value = 1.032240891231112
$ws.Range("L3").Value = 1.0531119374238
$ws.Range("M3").Value = 1.057706964228979
$ws.Range("N3").Value = 1.01424349075887
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.026641392387307
$ws.Range("D4").Value = 1.030117654949132
$ws.Range("E4").Value = 1.05128959477715
$ws.Range("F4").Value = 1.055923953303652
$ws.Range("I4").Value = 1.031570219809092
$ws.Range("J4").Value = 1.031217542711388
$ws.Range("K4").Value = 1.032616587082912
$ws.Range("L4").Value = 1.053735498795974
$ws.Range("M4").Value = 1.058358555012205
$ws.Range("N4").Value = 1.014402049341685
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.026937920556542
$ws.Range("D5").Value = 1.030328001622628
$ws.Range("E5").Value = 1.051603471821205
$ws.Range("F5").Value = 1.056249348976793
$ws.Range("I5").Value = 1.031620510913578
$ws.Range("J5").Value = 1.03141534499957
$ws.Range("K5").Value = 1.03277430305574
$ws.Range("L5").Value = 1.053997677326133
$ws.Range("M5").Value = 1.058632478261946
$ws.Range("N5").Value = 1.014468625335148
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.026987706868335
$ws.Range("D6").Value = 1.030363315644378
$ws.Range("E6").Value = 1.051656183284893
$ws.Range("F6").Value = 1.056303992083652
$ws.Range("I6").Value = 1.03162893466124
$ws.Range("J6").Value = 1.03144854735517
$ws.Range("K6").Value = 1.032800770973568
$ws.Range("L6").Value = 1.054041700224748
$ws.Range("M6").Value = 1.058678470872965
$ws.Range("N6").Value = 1.014479798926092
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.026645354750644
$ws.Range("D7").Value = 1.03012046588688
$ws.Range("E7").Value = 1.051293788139428
$ws.Range("F7").Value = 1.055928300738294
$ws.Range("I7").Value = 1.031570893164809
$ws.Range("J7").Value = 1.03122018639321
$ws.Range("K7").Value = 1.0326186953818
$ws.Range("L7").Value = 1.053739001906704
$ws.Range("M7").Value = 1.058362215209081
$ws.Range("N7").Value = 1.014402939255742
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025214073449095
$ws.Range("D8").Value = 1.029104800927504
$ws.Range("E8").Value = 1.049780497063738
$ws.Range("F8").Value = 1.05435908805996
$ws.Range("I8").Value = 1.031325344791606
$ws.Range("J8").Value = 1.030264295996724
$ws.Range("K8").Value = 1.031855726572744
$ws.Range("L8").Value = 1.052474053854997
$ws.Range("M8").Value = 1.057040269982601
$ws.Range("N8").Value = 1.014080982251343
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.022690768750859
$ws.Range("D9").Value = 1.027312855915029
$ws.Range("E9").Value = 1.047119352852269
$ws.Range("F9").Value = 1.051598095560749
$ws.Range("I9").Value = 1.03088156523723
$ws.Range("J9").Value = 1.028574684454535
$ws.Range("K9").Value = 1.03050403515171
$ws.Range("L9").Value = 1.050246079475803
$ws.Range("M9").Value = 1.054710616048975
$ws.Range("N9").Value = 1.013511030156329
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021007703297526
$ws.Range("D10").Value = 1.026116753261636
$ws.Range("E10").Value = 1.045348891973186
$ws.Range("F10").Value = 1.049760186020431
$ws.Range("I10").Value = 1.030578259823878
$ws.Range("J10").Value = 1.027444758750736
$ws.Range("K10").Value = 1.029598039184417
$ws.Range("L10").Value = 1.048761417262984
$ws.Range("M10").Value = 1.053157330881377
$ws.Range("N10").Value = 1.01312929890693
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020278700241392
$ws.Range("D11").Value = 1.025598480791448
$ws.Range("E11").Value = 1.044583115131548
$ws.Range("F11").Value = 1.048964992712425
$ws.Range("I11").Value = 1.030445158346531
$ws.Range("J11").Value = 1.026954649931512
$ws.Range("K11").Value = 1.029204577725751
$ws.Range("L11").Value = 1.048118688598167
$ws.Range("M11").Value = 1.052484689068702
$ws.Range("N11").Value = 1.012963586566835
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02000788110875
$ws.Range("D12").Value = 1.025405918517907
$ws.Range("E12").Value = 1.044298797425227
$ws.Range("F12").Value = 1.04866971666519
$ws.Range("I12").Value = 1.03039545287183
$ws.Range("J12").Value = 1.026772474510435
$ws.Range("K12").Value = 1.029058254319529
$ws.Range("L12").Value = 1.047879970731402
$ws.Range("M12").Value = 1.052234830207282
$ws.Range("N12").Value = 1.012901970431439
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020065974332254
$ws.Range("D13").Value = 1.025447226167078
$ws.Range("E13").Value = 1.044359778856025
$ws.Range("F13").Value = 1.048733050135806
$ws.Range("I13").Value = 1.030406126886593
$ws.Range("J13").Value = 1.026811557519273
$ws.Range("K13").Value = 1.029089649076679
$ws.Range("L13").Value = 1.047931175624559
$ws.Range("M13").Value = 1.052288426236004
$ws.Range("N13").Value = 1.012915190162364
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020256314967052
$ws.Range("D14").Value = 1.025582564610423
$ws.Range("E14").Value = 1.044559610778254
$ws.Range("F14").Value = 1.048940583187592
$ws.Range("I14").Value = 1.03044105509588
$ws.Range("J14").Value = 1.026939593851504
$ws.Range("K14").Value = 1.029192486143872
$ws.Range("L14").Value = 1.048098955698067
$ws.Range("M14").Value = 1.052464035861028
$ws.Range("N14").Value = 1.012958494647381
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02037358544874
$ws.Range("D15").Value = 1.02566594409827
$ws.Range("E15").Value = 1.044682750476491
$ws.Range("F15").Value = 1.049068463595229
$ws.Range("I15").Value = 1.030462540315823
$ws.Range("J15").Value = 1.027018464394911
$ws.Range("K15").Value = 1.029255824354103
$ws.Range("L15").Value = 1.048202333187407
$ws.Range("M15").Value = 1.052572233429953
$ws.Range("N15").Value = 1.012985167594046
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021056079995328
$ws.Range("D16").Value = 1.026151141888406
$ws.Range("E16").Value = 1.045399731769353
$ws.Range("F16").Value = 1.04981297364241
$ws.Range("I16").Value = 1.030587056075726
$ws.Range("J16").Value = 1.027477267848113
$ws.Range("K16").Value = 1.029624127490063
$ws.Range("L16").Value = 1.048804075923125
$ws.Range("M16").Value = 1.053201970584404
$ws.Range("N16").Value = 1.013140287830734
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021484129817082
$ws.Range("D17").Value = 1.026455399441849
$ws.Range("E17").Value = 1.045849700589685
$ws.Range("F17").Value = 1.050280154033728
$ws.Range("I17").Value = 1.030664688114633
$ws.Range("J17").Value = 1.027764836976872
$ws.Range("K17").Value = 1.029854843860648
$ws.Range("L17").Value = 1.049181570099776
$ws.Range("M17").Value = 1.053596971634011
$ws.Range("N17").Value = 1.013237478114739
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.021733782444803
$ws.Range("D18").Value = 1.026632833743363
$ws.Range("E18").Value = 1.046112241326173
$ws.Range("F18").Value = 1.050552713660177
$ws.Range("I18").Value = 1.030709798960984
$ws.Range("J18").Value = 1.027932489826148
$ws.Range("K18").Value = 1.029989304977304
$ws.Range("L18").Value = 1.049401769617911
$ws.Range("M18").Value = 1.053827363420561
$ws.Range("N18").Value = 1.013294127025279
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.021818903913292
$ws.Range("D19").Value = 1.026693328487812
$ws.Range("E19").Value = 1.046201774791652
$ws.Range("F19").Value = 1.050645659863023
$ws.Range("I19").Value = 1.030725151655735
$ws.Range("J19").Value = 1.027989641329164
$ws.Range("K19").Value = 1.030035133777724
$ws.Range("L19").Value = 1.049476854287191
$ws.Range("M19").Value = 1.053905920135153
$ws.Range("N19").Value = 1.013313435965879
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021438206363932
$ws.Range("D20").Value = 1.026422758999875
$ws.Range("E20").Value = 1.045801414766732
$ws.Range("F20").Value = 1.050230023675391
$ws.Range("I20").Value = 1.030656376578369
$ws.Range("J20").Value = 1.027733991949275
$ws.Range("K20").Value = 1.02983010174649
$ws.Range("L20").Value = 1.049141067151372
$ws.Range("M20").Value = 1.053554592349383
$ws.Range("N20").Value = 1.013227054710503
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.020200265368292
$ws.Range("D21").Value = 1.025542712268884
$ws.Range("E21").Value = 1.044500761789784
$ws.Range("F21").Value = 1.048879467285765
$ws.Range("I21").Value = 1.03043077694234
$ws.Range("J21").Value = 1.026901893857414
$ws.Range("K21").Value = 1.029162208003453
$ws.Range("L21").Value = 1.048049548097965
$ws.Range("M21").Value = 1.052412323479312
$ws.Range("N21").Value = 1.012945744300872
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.019421719994123
$ws.Range("D22").Value = 1.024989086510683
$ws.Range("E22").Value = 1.043683717639408
$ws.Range("F22").Value = 1.04803086310549
$ws.Range("I22").Value = 1.030287396470497
$ws.Range("J22").Value = 1.026377984691105
$ws.Range("K22").Value = 1.028741268171453
$ws.Range("L22").Value = 1.047363383070207
$ws.Range("M22").Value = 1.051694076907205
$ws.Range("N22").Value = 1.012768507627221
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.01983446124368
$ws.Range("D23").Value = 1.025282602920608
$ws.Range("E23").Value = 1.044116779514669
$ws.Range("F23").Value = 1.048480673000328
$ws.Range("I23").Value = 1.030363550864935
$ws.Range("J23").Value = 1.026655788775469
$ws.Range("K23").Value = 1.02896451198545
$ws.Range("L23").Value = 1.04772712138734
$ws.Range("M23").Value = 1.052074838677511
$ws.Range("N23").Value = 1.012862498823276
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021458957261993
$ws.Range("D24").Value = 1.026437507915495
$ws.Range("E24").Value = 1.045823232799431
$ws.Range("F24").Value = 1.050252675237013
$ws.Range("I24").Value = 1.030660132730929
$ws.Range("J24").Value = 1.027747929740708
$ws.Range("K24").Value = 1.029841281988286
$ws.Range("L24").Value = 1.049159368648165
$ws.Range("M24").Value = 1.053573741741039
$ws.Range("N24").Value = 1.013231764723403
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.023343251621929
$ws.Range("D25").Value = 1.027776377343587
$ws.Range("E25").Value = 1.04780667866926
$ws.Range("F25").Value = 1.052311389904202
$ws.Range("I25").Value = 1.030997606801023
$ws.Range("J25").Value = 1.02901210852059
$ws.Range("K25").Value = 1.03085433795131
$ws.Range("L25").Value = 1.050821946373695
$ws.Range("M25").Value = 1.055312916950846
$ws.Range("N25").Value = 1.013658687215156
